$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.6
$ws.Range("H2").Value = 5
$ws.Range("L2").Value = 1.33
$ws.Range("P2").Value = 1.88
$ws.Range("T2").Value = 1.89
$ws.Range("U2").Value = 1.87
$ws.Range("V2").Value = 1.15
$ws.Range("Z2").Value = 60
$ws.Range("AI2").Value = 110
$ws.Range("F3").Value = 2.18
$ws.Range("G3").Value = 2.38
$ws.Range("H3").Value = 3.3
$ws.Range("K3").Value = 3.85
$ws.Range("L3").Value = 1.41
$ws.Range("O3").Value = 1.33
$ws.Range("P3").Value = 1.8
$ws.Range("Q3").Value = 1.98
$ws.Range("S3").Value = 3.5
$ws.Range("W3").Value = 1.72
$ws.Range("AO3").Value = 50
$ws.Range("L4").Value = 1.48
$ws.Range("P4").Value = 1.67
$ws.Range("G5").Value = 4.1
$ws.Range("H5").Value = 2.16
$ws.Range("I5").Value = 2.4
$ws.Range("L5").Value = 1.45
$ws.Range("N5").Value = 3.1
$ws.Range("Y5").Value = 9
$ws.Range("AH5").Value = 980
$ws.Range("J6").Value = 3.15
$ws.Range("Q6").Value = 2.24
$ws.Range("V6").Value = 1.68
$ws.Range("Y8").Value = 18.5
$ws.Range("AE8").Value = 65
$ws.Range("AI8").Value = 70
$ws.Range("AL8").Value = 980
$ws.Range("G9").Value = 3.2
$ws.Range("H9").Value = 2.38
$ws.Range("I9").Value = 2.66
$ws.Range("J9").Value = 3.45
$ws.Range("K9").Value = 3.95
$ws.Range("O9").Value = 1.25
$ws.Range("P9").Value = 2.06
$ws.Range("Q9").Value = 1.75
$ws.Range("R9").Value = 1.42
$ws.Range("S9").Value = 2.88
$ws.Range("V9").Value = 1.6
$ws.Range("W9").Value = 1.46
$ws.Range("X9").Value = 18.5
$ws.Range("Y9").Value = 13
$ws.Range("Z9").Value = 18.5
$ws.Range("AB9").Value = 14.5
$ws.Range("AE9").Value = 980
$ws.Range("AF9").Value = 23
$ws.Range("AH9").Value = 16.5
$ws.Range("AI9").Value = 980
$ws.Range("AJ9").Value = 50
$ws.Range("AK9").Value = 34
$ws.Range("AN9").Value = 26
$ws.Range("AO9").Value = 18.5
$ws.Range("N10").Value = 3.5
$ws.Range("P10").Value = 1.85
$ws.Range("Q10").Value = 2.04
$ws.Range("S10").Value = 3.7
$ws.Range("T10").Value = 1.92
$ws.Range("AB10").Value = 8
$ws.Range("AE10").Value = 80
$ws.Range("AN10").Value = 13.5
$ws.Range("F11").Value = 2.6
$ws.Range("G11").Value = 2.74
$ws.Range("H11").Value = 2.9
$ws.Range("T12").Value = 2.04
$ws.Range("F13").Value = 2.32
$ws.Range("H13").Value = 3.15
$ws.Range("K13").Value = 3.85
$ws.Range("L13").Value = 1.33
$ws.Range("O13").Value = 1.24
$ws.Range("P13").Value = 2.22
$ws.Range("Q13").Value = 1.72
$ws.Range("R13").Value = 1.49
$ws.Range("S13").Value = 2.78
$ws.Range("U13").Value = 2.38
$ws.Range("V13").Value = 1.44
$ws.Range("J15").Value = 3.3
$ws.Range("O15").Value = 1.35
$ws.Range("S17").Value = 2.62
$ws.Range("K18").Value = 5.4
$ws.Range("L19").Value = 1.3
$ws.Range("F20").Value = 2.82
$ws.Range("I20").Value = 3.15
$ws.Range("J20").Value = 3
$ws.Range("K20").Value = 3.1
$ws.Range("N20").Value = 2.66
$ws.Range("P20").Value = 1.56
$ws.Range("S20").Value = 5.3
$ws.Range("T20").Value = 2.06
$ws.Range("V20").Value = 1.46
$ws.Range("J21").Value = 2.96
$ws.Range("N21").Value = 2.76
$ws.Range("O21").Value = 1.46
$ws.Range("P21").Value = 1.6
$ws.Range("Q21").Value = 2.34
$ws.Range("R21").Value = 1.22
$ws.Range("U21").Value = 1.84
$ws.Range("AH21").Value = 22
$ws.Range("AM21").Value = 180
$ws.Range("I22").Value = 2.18
$ws.Range("J22").Value = 3.25
$ws.Range("N22").Value = 3.1
$ws.Range("Q22").Value = 2.2
$ws.Range("T22").Value = 1.9
$ws.Range("V22").Value = 1.84
$ws.Range("W22").Value = 1.29
$ws.Range("AA22").Value = 28
$ws.Range("AD22").Value = 11
$ws.Range("AE22").Value = 27
$ws.Range("AG22").Value = 18
$ws.Range("AL22").Value = 75
$ws.Range("AO22").Value = 21
$ws.Range("F23").Value = 2.16
$ws.Range("G23").Value = 2.18
$ws.Range("I23").Value = 4.4
$ws.Range("N23").Value = 3.1
$ws.Range("Q23").Value = 2.4
$ws.Range("U23").Value = 1.91
$ws.Range("V23").Value = 1.29
$ws.Range("W23").Value = 1.84
$ws.Range("F24").Value = 2.28
$ws.Range("K24").Value = 3
$ws.Range("N24").Value = 2.16
$ws.Range("AB24").Value = 6.4
$ws.Range("AC24").Value = 7.4
$ws.Range("R25").Value = 1.47
$ws.Range("S25").Value = 3.05
$ws.Range("AE25").Value = 60
$ws.Range("AM25").Value = 90
$ws.Range("F26").Value = 2.08
$ws.Range("G26").Value = 2.1
$ws.Range("L26").Value = 1.49
$ws.Range("O26").Value = 1.42
$ws.Range("T26").Value = 1.98
$ws.Range("W26").Value = 1.9
$ws.Range("F27").Value = 1.64
$ws.Range("P27").Value = 1.74
$ws.Range("Q27").Value = 2.16
$ws.Range("N28").Value = 1.3
$ws.Range("P28").Value = 1.3
$ws.Range("I29").Value = 3
$ws.Range("K29").Value = 3.3
$ws.Range("Q29").Value = 2.3
$ws.Range("V29").Value = 1.5
$ws.Range("G30").Value = 2.42
$ws.Range("I30").Value = 4.3
$ws.Range("M30").Value = 1.07
$ws.Range("N30").Value = 3.55
$ws.Range("O30").Value = 1.33
$ws.Range("P30").Value = 1.86
$ws.Range("Q30").Value = 1.98
$ws.Range("AG30").Value = 980
$ws.Range("J31").Value = 3.25
$ws.Range("K31").Value = 3.55
$ws.Range("O31").Value = 1.39
$ws.Range("V31").Value = 1.35
$ws.Range("G32").Value = 1.8
$ws.Range("I32").Value = 7.2
$ws.Range("J32").Value = 3.5
$ws.Range("K32").Value = 3.75
$ws.Range("N32").Value = 2.58
$ws.Range("Q32").Value = 2.58
